# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gains a new (blank) column between the
# existing "In Advance" (M) and "Late" (N) columns, shifting "Late",
# "Paid" and "Outstanding" one column to the right. The "Repayment
# schedule" tab also becomes the active/selected sheet (it was
# "Transactions" before).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (14) - shifts N,O,P -> O,P,Q
$ws.Columns("N:N").Insert()

# New column takes the same width as the column to its left ("In Advance").
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, and move the selection
# to where the author last clicked.
$ws.Activate()
$ws.Range("J18").Select()
